# repull data, push all data, mean calculation
# Update the dSF (column F) values on Sheet1 to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = 4
    5  = 2
    6  = -2
    8  = -4
    11 = -5
    14 = 1
    15 = -4
    17 = -2
    20 = -3
    23 = -3
    26 = -1
    28 = -5
    29 = 1
    32 = -11
    33 = -2
    34 = -2
    35 = -1
    37 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
